$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 193.95454
$ws.Range("I33").Value = 208.78947
$ws.Range("K33").Value = 208.78947
$ws.Range("M33").Value = 20.21053000000001
$ws.Range("H87").Value = 26500
$ws.Range("J87").Value = 26500
$ws.Range("L87").Value = 26500
$ws.Range("N87").Value = -28996
$ws.Range("H90").Value = 26500
$ws.Range("J90").Value = 26500
$ws.Range("L90").Value = 79500
$ws.Range("N90").Value = -91980
$ws.Range("H125").Value = 1505.1578
$ws.Range("I125").Value = 690
$ws.Range("J125").Value = 2626
$ws.Range("K125").Value = 6210
$ws.Range("L125").Value = 23634
$ws.Range("M125").Value = -3750
$ws.Range("N125").Value = -28554
$ws.Range("H138").Value = 3710.7322
$ws.Range("I138").Value = 1472.2727
$ws.Range("K138").Value = 4416.8181
$ws.Range("M138").Value = 723.1818999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3341.9666
$ws.Range("I61").Value = 2625.7144
$ws.Range("J61").Value = 3968.6875
$ws.Range("K61").Value = 2625.7144
$ws.Range("L61").Value = 3968.6875
$ws.Range("M61").Value = -2413.7144
$ws.Range("N61").Value = -4392.6875
$ws.Range("H74").Value = 2391.182
$ws.Range("I74").Value = 2114.3103
$ws.Range("J74").Value = 4398.5
$ws.Range("K74").Value = 2114.3103
$ws.Range("L74").Value = 4398.5
$ws.Range("M74").Value = -1240.3103
$ws.Range("N74").Value = -6146.5
$ws.Range("H77").Value = 2391.182
$ws.Range("I77").Value = 2114.3103
$ws.Range("J77").Value = 4398.5
$ws.Range("K77").Value = 10571.5515
$ws.Range("L77").Value = 21992.5
$ws.Range("M77").Value = -6203.551500000001
$ws.Range("N77").Value = -30728.5
$ws.Range("H110").Value = 2559.8
$ws.Range("I110").Value = 1266.3334
$ws.Range("K110").Value = 1266.3334
$ws.Range("M110").Value = 778.6666
$ws.Range("H122").Value = 11145.55
$ws.Range("I122").Value = 11205.842
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 33617.526
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -31167.526
$ws.Range("N122").Value = -34900
$ws.Range("H136").Value = 3341.9666
$ws.Range("I136").Value = 2625.7144
$ws.Range("J136").Value = 3968.6875
$ws.Range("K136").Value = 7877.1432
$ws.Range("L136").Value = 11906.0625
$ws.Range("M136").Value = -5327.1432
$ws.Range("N136").Value = -17006.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9999.6
$ws.Range("J4").Value = 9999.6
$ws.Range("L4").Value = 9999.6
$ws.Range("N4").Value = -10223.6
$ws.Range("H31").Value = 4424.391
$ws.Range("I31").Value = 2099.7856
$ws.Range("J31").Value = 8040.4443
$ws.Range("K31").Value = 2099.7856
$ws.Range("L31").Value = 8040.4443
$ws.Range("M31").Value = -1804.7856
$ws.Range("N31").Value = -8630.444299999999
$ws.Range("H34").Value = 4424.391
$ws.Range("I34").Value = 2099.7856
$ws.Range("J34").Value = 8040.4443
$ws.Range("K34").Value = 2099.7856
$ws.Range("L34").Value = 8040.4443
$ws.Range("M34").Value = -1897.7856
$ws.Range("N34").Value = -8444.444299999999
$ws.Range("H50").Value = 15205.5
$ws.Range("J50").Value = 15205.5
$ws.Range("L50").Value = 15205.5
$ws.Range("N50").Value = -16455.5
$ws.Range("H51").Value = 17641.666
$ws.Range("J51").Value = 18962.5
$ws.Range("L51").Value = 18962.5
$ws.Range("N51").Value = -20434.5
$ws.Range("H61").Value = 17641.666
$ws.Range("J61").Value = 18962.5
$ws.Range("L61").Value = 18962.5
$ws.Range("N61").Value = -19658.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 139
$ws.Range("I4").Value = 93.333336
$ws.Range("J4").Value = 207.5
$ws.Range("K4").Value = 280.000008
$ws.Range("L4").Value = 622.5
$ws.Range("M4").Value = -168.000008
$ws.Range("N4").Value = -846.5
$ws.Range("H5").Value = 2044.625
$ws.Range("J5").Value = 885.4286
$ws.Range("L5").Value = 2656.2858
$ws.Range("N5").Value = -2880.2858
$ws.Range("H68").Value = 1111721.4
$ws.Range("I68").Value = 2000528
$ws.Range("J68").Value = 713
$ws.Range("K68").Value = 6001584
$ws.Range("L68").Value = 2139
$ws.Range("M68").Value = -6000773
$ws.Range("N68").Value = -3761
$ws.Range("H71").Value = 1111721.4
$ws.Range("I71").Value = 2000528
$ws.Range("J71").Value = 713
$ws.Range("K71").Value = 18004752
$ws.Range("L71").Value = 6417
$ws.Range("M71").Value = -18000696
$ws.Range("N71").Value = -14529
$ws.Range("H109").Value = 10489.77
$ws.Range("I109").Value = 34675.668
$ws.Range("K109").Value = 104027.004
$ws.Range("M109").Value = -102987.004
$ws.Range("H112").Value = 6235
$ws.Range("J112").Value = 6235
$ws.Range("L112").Value = 18705
$ws.Range("N112").Value = -20921
$ws.Range("H128").Value = 170000
$ws.Range("I128").Value = 170000
$ws.Range("K128").Value = 510000
$ws.Range("M128").Value = -505020
$ws.Range("H135").Value = 2044.625
$ws.Range("J135").Value = 885.4286
$ws.Range("L135").Value = 7968.8574
$ws.Range("N135").Value = -13038.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1940.3077
$ws.Range("I126").Value = 1685.3334
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 5056.0002
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -2586.0002
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9145.883
$ws.Range("H40").Value = 3485.5625
$ws.Range("I40").Value = 3461.5
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 3461.5
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -3325.5
$ws.Range("N40").Value = -3772
$ws.Range("H132").Value = 3246.795
$ws.Range("I132").Value = 3183.3447
$ws.Range("J132").Value = 3430.8
$ws.Range("K132").Value = 9550.034100000001
$ws.Range("L132").Value = 10292.4
$ws.Range("M132").Value = -7020.034100000001
$ws.Range("N132").Value = -15352.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H62").Value = 4550
$ws.Range("H65").Value = 4550
$ws.Range("H122").Value = 12502649
$ws.Range("J122").Value = 7750
$ws.Range("L122").Value = 23250
$ws.Range("N122").Value = -28150
